# Timers.xlsx - "Updated Changes TIMER Driver"
# The prescaler/duty-cycle inputs on Sheet1 were retuned, which ripples
# through the dependent ARR formula, and the active selection moved to C25.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# TIM_Prescaler: 100 -> 50
$ws1.Range("C24").Value = 50

# Duty cycle (%): 25 -> 50
$ws1.Range("C25").Value = 50

# Leave the cursor on C25, matching the saved selection in the workbook.
$ws1.Activate()
$ws1.Range("C25").Select() | Out-Null
